# Daily attendance processing - 2026-01-10 15:55:26
# Normalize the "Recorded By" column (G): move any entry that is the
# automated "system" recorder (case-insensitive match on "system") to the
# end of the comma-separated list, preserving the relative order of the
# remaining (human) recorders and of the system entries themselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ", "

    $nonSystem = @()
    $systemLike = @()
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $systemLike += $p
        } else {
            $nonSystem += $p
        }
    }

    $ordered = $nonSystem + $systemLike
    $newText = [string]::Join(", ", $ordered)

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
